$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Text corrections on the "CAM" sheet (renumbering two list items)
# ---------------------------------------------------------------------------
$camSheet = $wb.Worksheets.Item("CAM")
$camSheet.Range("A4").Value = "2. Camera BasicAuth"
$camSheet.Range("A5").Value = "3. Server-Hosted Camera UI"

# ---------------------------------------------------------------------------
# 2) Selections on each sheet (mirrors the saved view/selection state)
# ---------------------------------------------------------------------------
$sitegroundSheet = $wb.Worksheets.Item("Siteground")
[void]$sitegroundSheet.Range("A5:A6").Select()
[void]$sitegroundSheet.Range("L16").Select()

$serverSheet = $wb.Worksheets.Item("SERVER")
[void]$serverSheet.Range("A5:A6").Select()
[void]$serverSheet.Range("A4").Select()

[void]$camSheet.Range("A5:A6").Select()

# ---------------------------------------------------------------------------
# 3) Active sheet becomes "CAM" (the third tab, index 2)
# ---------------------------------------------------------------------------
[void]$camSheet.Activate()
